# REPORTGEN-1102: part 1, added and removed counts missing when no previous snapshot selected
#
# The RepGen table placeholders that drive the "Quality standards evolution" and
# "Rules list statistics ratio" tables need an extra EVOLUTION=true argument so that
# the report generator includes added/removed vulnerability counts even when no
# previous snapshot was selected.

$wb = $excel.ActiveWorkbook

# Map each worksheet to the cell holding the RepGen placeholder that must be updated.
$targets = @(
    @{ Sheet = "Summary";     Cell = "B14" },
    @{ Sheet = "API1-2019";   Cell = "A3" },
    @{ Sheet = "API2-2019";   Cell = "A3" },
    @{ Sheet = "API3-2019";   Cell = "A3" },
    @{ Sheet = "API4-2019";   Cell = "A3" },
    @{ Sheet = "API5-2019";   Cell = "A3" },
    @{ Sheet = "API6-2019";   Cell = "A3" },
    @{ Sheet = "API7-2019";   Cell = "A3" },
    @{ Sheet = "API8-2019";   Cell = "A3" },
    @{ Sheet = "API9-2019";   Cell = "A3" },
    @{ Sheet = "API10-2019";  Cell = "A3" }
)

foreach ($target in $targets) {
    $ws = $wb.Worksheets.Item($target.Sheet)
    $range = $ws.Range($target.Cell)
    $current = $range.Value()
    if ($current -notlike "*,EVOLUTION=true") {
        $range.Value = $current + ",EVOLUTION=true"
    }
}
